$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (HOUR=17): update existing B18 and D18
$ws.Range("B18").Value = 58430
$ws.Range("D18").Value = 41665.5135

# Row 19 (HOUR=18): add B19, update D19
$ws.Range("B19").Value = 56357
$ws.Range("D19").Value = 39727.5055

# Row 20 (HOUR=19): add B20, update D20
$ws.Range("B20").Value = 55191
$ws.Range("D20").Value = 39165.7675

# Row 21 (HOUR=20): add B21, update D21
$ws.Range("B21").Value = 52357
$ws.Range("D21").Value = 37817.495

# Row 22 (HOUR=21): add B22, update D22
$ws.Range("B22").Value = 50614
$ws.Range("D22").Value = 38010.227

# Row 23 (HOUR=22): add B23, update D23
$ws.Range("B23").Value = 47951
$ws.Range("D23").Value = 38285.2065
